$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row from the sheet's used range (falls back to 92).
$lastRow = $ws.UsedRange.Rows.Count
if (-not $lastRow -or $lastRow -lt 92) {
    $lastRow = 92
}

# 1) Column C ("Förändrad") changed from 45184 to 45186 for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -eq 45184) {
        $cell.Value = 45186
    }
}

# 2) HYPERLINK() formulas in columns S,T,U,V,W,X,Y gain a second argument:
#    the row's "Beteckning" (column A) text, used as the link's friendly label.
$hyperlinkCols = @("S","T","U","V","W","X","Y")

for ($r = 2; $r -le $lastRow; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    if (-not $beteckning) {
        continue
    }
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $r)
        $oldFormula = $cell.Formula
        if ($oldFormula -and $oldFormula.StartsWith("=HYPERLINK(")) {
            $argsStart = $oldFormula.IndexOf("(") + 1
            $argsEnd = $oldFormula.LastIndexOf(")")
            $argsText = $oldFormula.Substring($argsStart, $argsEnd - $argsStart)
            # Only rewrite formulas that don't already carry a second argument.
            $quoteStart = $argsText.IndexOf('"')
            $quoteEnd = $argsText.IndexOf('"', $quoteStart + 1)
            $url = $argsText.Substring($quoteStart + 1, $quoteEnd - $quoteStart - 1)
            $rest = $argsText.Substring($quoteEnd + 1).Trim()
            if (-not $rest) {
                $newFormula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
